$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 61

$ws.Range("A60:V60").Copy()
$ws.Range("A61:V61").PasteSpecial(-4122)

$ws.Cells.Item($row, 1).Value = 60
$ws.Cells.Item($row, 2).Value = "iran"
$ws.Cells.Item($row, 3).Value = "persian-gulf-pro-league"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45232.64583333334
$ws.Cells.Item($row, 6).Value = "Persepolis"
$ws.Cells.Item($row, 7).Value = 2
$ws.Cells.Item($row, 8).Value = "Sanat Naft"
$ws.Cells.Item($row, 9).Value = 2
$ws.Cells.Item($row, 10).Value = 1.3
$ws.Cells.Item($row, 11).Value = "01/11/2023 03:42"
$ws.Cells.Item($row, 12).Value = 1.2
$ws.Cells.Item($row, 13).Value = "02/11/2023 15:27"
$ws.Cells.Item($row, 14).Value = 4.44
$ws.Cells.Item($row, 15).Value = "01/11/2023 03:42"
$ws.Cells.Item($row, 16).Value = 5.77
$ws.Cells.Item($row, 17).Value = "02/11/2023 15:28"
$ws.Cells.Item($row, 18).Value = 8.529999999999999
$ws.Cells.Item($row, 19).Value = "01/11/2023 03:42"
$ws.Cells.Item($row, 20).Value = 15.87
$ws.Cells.Item($row, 21).Value = "02/11/2023 15:28"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/iran/persian-gulf-pro-league/persepolis-sanat-naft/U9y2Cmvc/"
